$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Rename the sheet tabs: "0~1".."9~10"  ->  "one".."ten"
# ------------------------------------------------------------------
$newNames = @("one","two","three","four","five","six","seven","eight","nine","ten")
for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i + 1).Name = $newNames[$i]
}

# ------------------------------------------------------------------
# 2) Clean up the header row text on every sheet (trim stray leading
#    spaces / slashes so the labels read like JSON field names).
# ------------------------------------------------------------------
$headers = @("TimeInterval", "Function", "HitCount", "Percentage", "Time(ms)", "Index", "TimeHit")

# Sheets 1-9 ("one".."nine") use columns A:G for the header row.
$cols = @("A", "B", "C", "D", "E", "F", "G")
for ($s = 1; $s -le 9; $s++) {
    $ws = $wb.Worksheets.Item($s)
    for ($c = 0; $c -lt $headers.Count; $c++) {
        $ws.Range($cols[$c] + "1").Value = $headers[$c]
    }
}

# Sheet 10 ("ten") uses columns O:U for the header row.
$cols10 = @("O", "P", "Q", "R", "S", "T", "U")
$ws10 = $wb.Worksheets.Item(10)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws10.Range($cols10[$c] + "1").Value = $headers[$c]
}

# ------------------------------------------------------------------
# 3) Update each sheet's selected cell / range (view state).
#    Sheet 1 ("one") is selected last so it ends up the active tab.
# ------------------------------------------------------------------
$wb.Worksheets.Item(10).Range("U3").Select()  | Out-Null   # ten
$wb.Worksheets.Item(2).Range("G4").Select()   | Out-Null   # two
$wb.Worksheets.Item(3).Range("E16").Select()  | Out-Null   # three
$wb.Worksheets.Item(4).Range("H5").Select()   | Out-Null   # four
$wb.Worksheets.Item(5).Range("F4").Select()   | Out-Null   # five
$wb.Worksheets.Item(6).Range("K4").Select()   | Out-Null   # six
$wb.Worksheets.Item(7).Range("H4").Select()   | Out-Null   # seven
$wb.Worksheets.Item(8).Range("H17").Select()  | Out-Null   # eight
$wb.Worksheets.Item(9).Range("E1").Select()   | Out-Null   # nine
$wb.Worksheets.Item(1).Range("A1").Select()   | Out-Null   # one (active tab)
